# Update PSSM score matrix (B2:K21) with recalculated values from supplemental figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 2.347406399582858, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, 4.321925899727339, 2.921125178424078, -19.2520093448428, 2.695884735718626, -19.2520093448428, 1.86804940477054, -19.2520093448428, 2.390984087137551, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 2.060012185346181, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(2.791345416539189, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, 2.846416361490208, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(3.708984015872135, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 1.515742718048111, -19.2520093448428, 2.239785536095939),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, 1.945440742951913, -19.2520093448428, 2.573162576227261, -19.2520093448428, -19.2520093448428, -19.2520093448428, 1.381118039820615),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, 1.679256948205758, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 2.278639406261707, 1.639656452390126),
    @(-19.2520093448428, -19.2520093448428, 1.382465838417955, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 2.085454130898911),
    @(-19.2520093448428, -19.2520093448428, -1.477318787227693, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 2.302808525636706, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -0.6616818010497065, -19.2520093448428, -19.2520093448428, -19.2520093448428, 0.4608954350008291, 0.8672775719928107, 1.270358587556165, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, -19.2520093448428, 0.3472079666209922, 0.8394996475471049, 1.354020779985607, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, 2.005257057250303, -19.2520093448428, -19.2520093448428, -19.2520093448428, 1.920277146636591, 2.145690875429304, -19.2520093448428, -19.2520093448428),
    @(-19.2520093448428, -19.2520093448428, 2.266809467388764, -19.2520093448428, 3.757072447022189, -19.2520093448428, 2.206469863536675, 2.004589239085011, -19.2520093448428, 2.409255985354795),
    @(-19.2520093448428, -19.2520093448428, -19.2520093448428, 2.52448352769484, -19.2520093448428, 3.304357545627962, 2.406181649691677, -19.2520093448428, -19.2520093448428, -19.2520093448428)
)

$arr = New-Object "object[,]" 20,10
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt $rows[$i].Count; $j++) {
        $arr[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("B2:K21").Value = $arr
